# "v1 e v2 a funcionar" — fill in the missing V2 measurement columns (F:H)
# on sheet V1 so the shared AVERAGE formulas in row 12 resolve instead of
# producing #DIV/0!, and extend the table with a new H column (16384/48).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row (row 1): new column H ---------------------------------
$ws.Range("H1").Value = 48

# --- Data rows 2-11: columns F, G, H -----------------------------------
$data = @{
    2  = @(10.797729,            8.6426719999999992,  7.5365580000000003)
    3  = @(10.795088,            8.5726130000000005,  7.9841420000000003)
    4  = @(10.769329000000001,   8.4481959999999994,  7.4888880000000002)
    5  = @(10.767128,            8.5706430000000005,  7.4482249999999999)
    6  = @(10.945178,            8.5990020000000005,  7.155856)
    7  = @(10.794335,            8.5672300000000003,  7.1304379999999998)
    8  = @(10.849694,            8.5589980000000008,  7.227595)
    9  = @(10.903616,            8.4926879999999993,  7.361364)
    10 = @(10.83319,             8.610989,             7.5399700000000003)
    11 = @(10.864084,            8.6026410000000002,  7.6326910000000003)
}

foreach ($row in $data.Keys | Sort-Object) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 6).Value = $vals[0]
    $ws.Cells.Item($row, 7).Value = $vals[1]
    $ws.Cells.Item($row, 8).Value = $vals[2]
}

# --- Row 12: extend the shared AVERAGE formula from F:G to F:H ---------
$ws.Range("F12:H12").Formula = "=AVERAGE(F2:F11)"

# --- Selection moves to H16 (matches the post-edit saved view) ---------
$ws.Range("H16").Select()
